$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2870.4333
$ws.Range("I98").Value = 2928.0344
$ws.Range("J98").Value = 1200
$ws.Range("K98").Value = 2928.0344
$ws.Range("L98").Value = 1200
$ws.Range("M98").Value = -1430.0344
$ws.Range("N98").Value = -4196

$ws.Range("H106").Value = 764.2632
$ws.Range("I106").Value = 677.0833
$ws.Range("J106").Value = 913.7143
$ws.Range("K106").Value = 677.0833
$ws.Range("L106").Value = 913.7143
$ws.Range("M106").Value = -46.08330000000001
$ws.Range("N106").Value = -2175.7143

$ws.Range("H116").Value = 1656.6364
$ws.Range("I116").Value = 1799.75
$ws.Range("J116").Value = 1574.8572
$ws.Range("K116").Value = 1799.75
$ws.Range("L116").Value = 1574.8572
$ws.Range("M116").Value = 1642.25
$ws.Range("N116").Value = -8458.8572

$ws.Range("H122").Value = 2870.4333
$ws.Range("I122").Value = 2928.0344
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 8784.1032
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -6334.1032
$ws.Range("N122").Value = -8500

$ws.Range("H137").Value = 19232066
$ws.Range("I137").Value = 914.25
$ws.Range("J137").Value = 62502156
$ws.Range("K137").Value = 2742.75
$ws.Range("L137").Value = 187506468
$ws.Range("M137").Value = -192.75
$ws.Range("N137").Value = -187511568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 796883
$ws.Range("I2").Value = 1726.7368
$ws.Range("K2").Value = 1726.7368
$ws.Range("M2").Value = -1613.7368

$ws.Range("H32").Value = 1448.52
$ws.Range("I32").Value = 1474.6123
$ws.Range("J32").Value = 170
$ws.Range("K32").Value = 1474.6123
$ws.Range("L32").Value = 170
$ws.Range("M32").Value = -1187.6123
$ws.Range("N32").Value = -744

$ws.Range("H61").Value = 4834729
$ws.Range("I61").Value = 6177161.5
$ws.Range("J61").Value = 1970.8
$ws.Range("K61").Value = 6177161.5
$ws.Range("L61").Value = 1970.8
$ws.Range("M61").Value = -6176949.5
$ws.Range("N61").Value = -2394.8

$ws.Range("H74").Value = 17243720
$ws.Range("I74").Value = 27778612
$ws.Range("J74").Value = 4803.727
$ws.Range("K74").Value = 27778612
$ws.Range("L74").Value = 4803.727
$ws.Range("M74").Value = -27777738
$ws.Range("N74").Value = -6551.727

$ws.Range("H77").Value = 17243720
$ws.Range("I77").Value = 27778612
$ws.Range("J77").Value = 4803.727
$ws.Range("K77").Value = 138893060
$ws.Range("L77").Value = 24018.635
$ws.Range("M77").Value = -138888692
$ws.Range("N77").Value = -32754.635

$ws.Range("H110").Value = 1385.625
$ws.Range("I110").Value = 1275
$ws.Range("J110").Value = 1570
$ws.Range("K110").Value = 1275
$ws.Range("L110").Value = 1570
$ws.Range("M110").Value = 770
$ws.Range("N110").Value = -5660

$ws.Range("H116").Value = 796883
$ws.Range("I116").Value = 1726.7368
$ws.Range("K116").Value = 1726.7368
$ws.Range("M116").Value = 567.2632000000001

$ws.Range("H132").Value = 1192436.6
$ws.Range("I132").Value = 1570709.8
$ws.Range("J132").Value = 3578.2856
$ws.Range("K132").Value = 4712129.4
$ws.Range("L132").Value = 10734.8568
$ws.Range("M132").Value = -4709599.4
$ws.Range("N132").Value = -15794.8568

$ws.Range("H136").Value = 4834729
$ws.Range("I136").Value = 6177161.5
$ws.Range("J136").Value = 1970.8
$ws.Range("K136").Value = 18531484.5
$ws.Range("L136").Value = 5912.4
$ws.Range("M136").Value = -18528934.5
$ws.Range("N136").Value = -11012.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 796883
$ws.Range("I3").Value = 1726.7368
$ws.Range("K3").Value = 1726.7368
$ws.Range("M3").Value = -1612.7368

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H94").Value = 550.48
$ws.Range("I94").Value = 298.1905
$ws.Range("J94").Value = 1875
$ws.Range("K94").Value = 298.1905
$ws.Range("L94").Value = 1875
$ws.Range("M94").Value = 152.8095
$ws.Range("N94").Value = -2777

$ws.Range("H134").Value = 10006730
$ws.Range("I134").Value = 15394140
$ws.Range("J134").Value = 1542
$ws.Range("K134").Value = 46182420
$ws.Range("L134").Value = 4626
$ws.Range("M134").Value = -46179885
$ws.Range("N134").Value = -9696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 13786.833
$ws.Range("J74").Value = 15905.7
$ws.Range("L74").Value = 15905.7
$ws.Range("N74").Value = -17653.7

$ws.Range("H77").Value = 13786.833
$ws.Range("J77").Value = 15905.7
$ws.Range("L77").Value = 47717.10000000001
$ws.Range("N77").Value = -56453.10000000001

$ws.Range("H132").Value = 2994.6775
$ws.Range("I132").Value = 2676.2917
$ws.Range("J132").Value = 4086.2856
$ws.Range("K132").Value = 8028.875100000001
$ws.Range("L132").Value = 12258.8568
$ws.Range("M132").Value = -5498.875100000001
$ws.Range("N132").Value = -17318.8568

$ws.Range("H134").Value = 1627.8055
$ws.Range("I134").Value = 1250.2333
$ws.Range("J134").Value = 3515.6667
$ws.Range("K134").Value = 3750.699900000001
$ws.Range("L134").Value = 10547.0001
$ws.Range("M134").Value = -1215.699900000001
$ws.Range("N134").Value = -15617.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1078395.1
$ws.Range("J64").Value = 2335150
$ws.Range("L64").Value = 7005450
$ws.Range("N64").Value = -7005990

$ws.Range("H67").Value = 1078395.1
$ws.Range("J67").Value = 2335150
$ws.Range("L67").Value = 7005450
$ws.Range("N67").Value = -7007322

$ws.Range("H70").Value = 25844.625
$ws.Range("I70").Value = 46444.25
$ws.Range("K70").Value = 139332.75
$ws.Range("M70").Value = -139017.75

$ws.Range("H73").Value = 25844.625
$ws.Range("I73").Value = 46444.25
$ws.Range("K73").Value = 139332.75
$ws.Range("M73").Value = -138240.75

$ws.Range("H113").Value = 2410.5925
$ws.Range("I113").Value = 561.2727
$ws.Range("J113").Value = 3682
$ws.Range("K113").Value = 1683.8181
$ws.Range("L113").Value = 11046
$ws.Range("M113").Value = 486.1819
$ws.Range("N113").Value = -15386

$ws.Range("H141").Value = 5452.1816
$ws.Range("I141").Value = 5452.1816
$ws.Range("K141").Value = 16356.5448
$ws.Range("M141").Value = -11176.5448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2270.889
$ws.Range("I132").Value = 1086.6666
$ws.Range("J132").Value = 2863
$ws.Range("K132").Value = 3259.9998
$ws.Range("L132").Value = 8589
$ws.Range("M132").Value = -729.9998000000001
$ws.Range("N132").Value = -13649

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1513.6923
$ws.Range("I61").Value = 1514.8334
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1514.8334
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -1312.8334
$ws.Range("N61").Value = -1904

$ws.Range("H82").Value = 1523.1818
$ws.Range("I82").Value = 1422.7778
$ws.Range("J82").Value = 1975
$ws.Range("K82").Value = 1422.7778
$ws.Range("L82").Value = 1975
$ws.Range("M82").Value = -1061.7778
$ws.Range("N82").Value = -2697

$ws.Range("H85").Value = 1523.1818
$ws.Range("I85").Value = 1422.7778
$ws.Range("J85").Value = 1975
$ws.Range("K85").Value = 1422.7778
$ws.Range("L85").Value = 1975
$ws.Range("M85").Value = -174.7778000000001
$ws.Range("N85").Value = -4471

$ws.Range("H113").Value = 1513.6923
$ws.Range("I113").Value = 1514.8334
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1514.8334
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 655.1666
$ws.Range("N113").Value = -5840

$ws.Range("H132").Value = 5427.75
$ws.Range("I132").Value = 6234.207
$ws.Range("K132").Value = 18702.621
$ws.Range("M132").Value = -16172.621

$ws.Range("H136").Value = 1718.909
$ws.Range("I136").Value = 1426
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 4278
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -1728
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 745.875
$ws.Range("I107").Value = 818.2
$ws.Range("J107").Value = 713
$ws.Range("K107").Value = 2454.6
$ws.Range("L107").Value = 2139
$ws.Range("M107").Value = -534.6000000000004
$ws.Range("N107").Value = -5979

$ws.Range("H113").Value = 373.0345
$ws.Range("I113").Value = 402.42105
$ws.Range("J113").Value = 317.2
$ws.Range("K113").Value = 1207.26315
$ws.Range("L113").Value = 951.5999999999999
$ws.Range("M113").Value = 962.73685
$ws.Range("N113").Value = -5291.6

$ws.Range("H122").Value = 1513.95
$ws.Range("I122").Value = 1511.6
$ws.Range("J122").Value = 1521
$ws.Range("K122").Value = 4534.799999999999
$ws.Range("L122").Value = 4563
$ws.Range("M122").Value = -2084.799999999999
$ws.Range("N122").Value = -9463
